$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.019624042182269
$ws.Cells.Item(2, 4).Value = 1.025740588271574
$ws.Cells.Item(2, 5).Value = 1.023240541836242
$ws.Cells.Item(2, 6).Value = 1.030978082333441
$ws.Cells.Item(2, 9).Value = 1.029812577073279
$ws.Cells.Item(2, 10).Value = 1.024825755802575
$ws.Cells.Item(2, 11).Value = 1.028565274147346
$ws.Cells.Item(2, 12).Value = 1.026072560940633
$ws.Cells.Item(2, 13).Value = 1.03378753074492
$ws.Cells.Item(2, 14).Value = 1.012348173130203

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.020590418791168
$ws.Cells.Item(3, 4).Value = 1.026466847552044
$ws.Cells.Item(3, 5).Value = 1.02415380733719
$ws.Cells.Item(3, 6).Value = 1.032174922221027
$ws.Cells.Item(3, 9).Value = 1.030019213324944
$ws.Cells.Item(3, 10).Value = 1.025428879622925
$ws.Cells.Item(3, 11).Value = 1.029099120466812
$ws.Cells.Item(3, 12).Value = 1.026792373047591
$ws.Cells.Item(3, 13).Value = 1.034791796450857
$ws.Cells.Item(3, 14).Value = 1.012548980625058

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.021215791608395
$ws.Cells.Item(4, 4).Value = 1.026936415182208
$ws.Cells.Item(4, 5).Value = 1.024745180147674
$ws.Cells.Item(4, 6).Value = 1.032949382534049
$ws.Cells.Item(4, 9).Value = 1.03015105958008
$ws.Cells.Item(4, 10).Value = 1.025818630847728
$ws.Cells.Item(4, 11).Value = 1.029443518935331
$ws.Cells.Item(4, 12).Value = 1.027257944345747
$ws.Cells.Item(4, 13).Value = 1.035441086489644
$ws.Cells.Item(4, 14).Value = 1.012678698337316

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.02147871304554
$ws.Cells.Item(5, 4).Value = 1.027133731746426
$ws.Cells.Item(5, 5).Value = 1.024993895284819
$ws.Cells.Item(5, 6).Value = 1.033274972181574
$ws.Cells.Item(5, 9).Value = 1.030206041744393
$ws.Cells.Item(5, 10).Value = 1.025982359662429
$ws.Cells.Item(5, 11).Value = 1.029588055447174
$ws.Cells.Item(5, 12).Value = 1.027453623225223
$ws.Cells.Item(5, 13).Value = 1.0357139197907
$ws.Cells.Item(5, 14).Value = 1.012733179237172

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.021522859566312
$ws.Cells.Item(6, 4).Value = 1.027166856780635
$ws.Cells.Item(6, 5).Value = 1.02503566161248
$ws.Cells.Item(6, 6).Value = 1.03332964051582
$ws.Cells.Item(6, 9).Value = 1.030215247325614
$ws.Cells.Item(6, 10).Value = 1.026009843248274
$ws.Cells.Item(6, 11).Value = 1.029612309166724
$ws.Cells.Item(6, 12).Value = 1.027486475780789
$ws.Cells.Item(6, 13).Value = 1.035759722192579
$ws.Cells.Item(6, 14).Value = 1.012742323739914

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.021219304720061
$ws.Cells.Item(7, 4).Value = 1.026939052089306
$ws.Cells.Item(7, 5).Value = 1.024748503091559
$ws.Cells.Item(7, 6).Value = 1.032953733052478
$ws.Cells.Item(7, 9).Value = 1.030151796007472
$ws.Cells.Item(7, 10).Value = 1.025820819082476
$ws.Cells.Item(7, 11).Value = 1.02944545121693
$ws.Cells.Item(7, 12).Value = 1.027260559203041
$ws.Cells.Item(7, 13).Value = 1.035444732605025
$ws.Cells.Item(7, 14).Value = 1.012679426520017

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.019950620575569
$ws.Cells.Item(8, 4).Value = 1.025986107196082
$ws.Cells.Item(8, 5).Value = 1.023549094956981
$ws.Cells.Item(8, 6).Value = 1.031382555174758
$ws.Cells.Item(8, 9).Value = 1.029882795887058
$ws.Cells.Item(8, 10).Value = 1.025029689694129
$ws.Cells.Item(8, 11).Value = 1.028745904101144
$ws.Cells.Item(8, 12).Value = 1.026315864916457
$ws.Cells.Item(8, 13).Value = 1.034127038439522
$ws.Cells.Item(8, 14).Value = 1.012416082017542

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.017715527063516
$ws.Cells.Item(9, 4).Value = 1.024304097700898
$ws.Cells.Item(9, 5).Value = 1.021438899730204
$ws.Cells.Item(9, 6).Value = 1.028614102999686
$ws.Cells.Item(9, 9).Value = 1.029394546826253
$ws.Cells.Item(9, 10).Value = 1.023631730505812
$ws.Cells.Item(9, 11).Value = 1.027505307085685
$ws.Cells.Item(9, 12).Value = 1.024649721003781
$ws.Cells.Item(9, 13).Value = 1.031800964277688
$ws.Cells.Item(9, 14).Value = 1.011950373029926

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.016225806302164
$ws.Cells.Item(10, 4).Value = 1.023180928444738
$ws.Cells.Item(10, 5).Value = 1.020034369760817
$ws.Cells.Item(10, 6).Value = 1.026768537883335
$ws.Cells.Item(10, 9).Value = 1.029059498896121
$ws.Cells.Item(10, 10).Value = 1.022697167950825
$ws.Cells.Item(10, 11).Value = 1.026672962102429
$ws.Cells.Item(10, 12).Value = 1.023537998589346
$ws.Cells.Item(10, 13).Value = 1.030247446906113
$ws.Cells.Item(10, 14).Value = 1.011638793627977

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015580823396982
$ws.Cells.Item(11, 4).Value = 1.022694160176556
$ws.Cells.Item(11, 5).Value = 1.0194267377241
$ws.Cells.Item(11, 6).Value = 1.025969394213933
$ws.Cells.Item(11, 9).Value = 1.028912159121052
$ws.Cells.Item(11, 10).Value = 1.022291882311905
$ws.Cells.Item(11, 11).Value = 1.026311301886118
$ws.Cells.Item(11, 12).Value = 1.023056388479543
$ws.Cells.Item(11, 13).Value = 1.029574085600305
$ws.Cells.Item(11, 14).Value = 1.011503615903732

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.015341259467849
$ws.Cells.Item(12, 4).Value = 1.022513289143979
$ws.Cells.Item(12, 5).Value = 1.019201117573757
$ws.Cells.Item(12, 6).Value = 1.025672555405762
$ws.Cells.Item(12, 9).Value = 1.028857091049554
$ws.Cells.Item(12, 10).Value = 1.02214124923888
$ws.Cells.Item(12, 11).Value = 1.026176777945998
$ws.Cells.Item(12, 12).Value = 1.022877463316601
$ws.Cells.Item(12, 13).Value = 1.029323866658952
$ws.Cells.Item(12, 14).Value = 1.011453365691941

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.015392646225909
$ws.Cells.Item(13, 4).Value = 1.02255208946142
$ws.Cells.Item(13, 5).Value = 1.019249510168321
$ws.Cells.Item(13, 6).Value = 1.02573622841118
$ws.Cells.Item(13, 9).Value = 1.028868918697109
$ws.Cells.Item(13, 10).Value = 1.022173564700343
$ws.Cells.Item(13, 11).Value = 1.026205642253053
$ws.Cells.Item(13, 12).Value = 1.022915844891559
$ws.Cells.Item(13, 13).Value = 1.029377544124168
$ws.Cells.Item(13, 14).Value = 1.011464146303907

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.015561020734395
$ws.Cells.Item(14, 4).Value = 1.022679210614115
$ws.Cells.Item(14, 5).Value = 1.019408086227758
$ws.Cells.Item(14, 6).Value = 1.025944857457092
$ws.Cells.Item(14, 9).Value = 1.028907614101486
$ws.Cells.Item(14, 10).Value = 1.022279432809025
$ws.Cells.Item(14, 11).Value = 1.026300185912077
$ws.Cells.Item(14, 12).Value = 1.023041599155574
$ws.Cells.Item(14, 13).Value = 1.029553404521063
$ws.Cells.Item(14, 14).Value = 1.011499463004108

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015664763342558
$ws.Cells.Item(15, 4).Value = 1.022757525743656
$ws.Cells.Item(15, 5).Value = 1.019505800978745
$ws.Cells.Item(15, 6).Value = 1.026073400500039
$ws.Cells.Item(15, 9).Value = 1.028931410636385
$ws.Cells.Item(15, 10).Value = 1.022344649470291
$ws.Cells.Item(15, 11).Value = 1.026358412584435
$ws.Cells.Item(15, 12).Value = 1.023119076050372
$ws.Cells.Item(15, 13).Value = 1.029661744302928
$ws.Cells.Item(15, 14).Value = 1.011521217600437

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.016268613296622
$ws.Cells.Item(16, 4).Value = 1.023213224687445
$ws.Cells.Item(16, 5).Value = 1.020074707691344
$ws.Cells.Item(16, 6).Value = 1.026821574293862
$ws.Cells.Item(16, 9).Value = 1.029069229723082
$ws.Cells.Item(16, 10).Value = 1.022724052514717
$ws.Cells.Item(16, 11).Value = 1.026696938017631
$ws.Cells.Item(16, 12).Value = 1.023569956716332
$ws.Cells.Item(16, 13).Value = 1.030292121332942
$ws.Cells.Item(16, 14).Value = 1.01164775942514

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.016647412785172
$ws.Cells.Item(17, 4).Value = 1.023498958396961
$ws.Cells.Item(17, 5).Value = 1.020431712429665
$ws.Cells.Item(17, 6).Value = 1.027290882765021
$ws.Cells.Item(17, 9).Value = 1.029155074652331
$ws.Cells.Item(17, 10).Value = 1.022961877870968
$ws.Cells.Item(17, 11).Value = 1.026908951924764
$ws.Cells.Item(17, 12).Value = 1.023852721746368
$ws.Cells.Item(17, 13).Value = 1.030687358165656
$ws.Cells.Item(17, 14).Value = 1.011727065783177

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.016868367451271
$ws.Cells.Item(18, 4).Value = 1.023665580634138
$ws.Cells.Item(18, 5).Value = 1.020639999197256
$ws.Cells.Item(18, 6).Value = 1.027564622350145
$ws.Cells.Item(18, 9).Value = 1.029204928295909
$ws.Cells.Item(18, 10).Value = 1.023100538127476
$ws.Cells.Item(18, 11).Value = 1.027032495475072
$ws.Cells.Item(18, 12).Value = 1.024017631880369
$ws.Cells.Item(18, 13).Value = 1.030917827750512
$ws.Cells.Item(18, 14).Value = 1.011773298571817

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.016943708525645
$ws.Cells.Item(19, 4).Value = 1.02372238746538
$ws.Cells.Item(19, 5).Value = 1.020711028397425
$ws.Cells.Item(19, 6).Value = 1.027657960486909
$ws.Cells.Item(19, 9).Value = 1.02922189007186
$ws.Cells.Item(19, 10).Value = 1.02314780763297
$ws.Cells.Item(19, 11).Value = 1.027074600167197
$ws.Cells.Item(19, 12).Value = 1.024073858244235
$ws.Cells.Item(19, 13).Value = 1.030996400856079
$ws.Cells.Item(19, 14).Value = 1.011789058467633

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.0166067704262
$ws.Cells.Item(20, 4).Value = 1.023468306150819
$ws.Cells.Item(20, 5).Value = 1.020393403823346
$ws.Cells.Item(20, 6).Value = 1.027240530425503
$ws.Cells.Item(20, 9).Value = 1.029145886876427
$ws.Cells.Item(20, 10).Value = 1.022936367587348
$ws.Cells.Item(20, 11).Value = 1.026886217307155
$ws.Cells.Item(20, 12).Value = 1.023822385999609
$ws.Cells.Item(20, 13).Value = 1.030644959760691
$ws.Cells.Item(20, 14).Value = 1.01171855957405

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015511438322332
$ws.Cells.Item(21, 4).Value = 1.022641778353314
$ws.Cells.Item(21, 5).Value = 1.019361387278377
$ws.Cells.Item(21, 6).Value = 1.025883421446035
$ws.Cells.Item(21, 9).Value = 1.028896228639405
$ws.Cells.Item(21, 10).Value = 1.02224825982545
$ws.Cells.Item(21, 11).Value = 1.026272350325755
$ws.Cells.Item(21, 12).Value = 1.023004568587134
$ws.Cells.Item(21, 13).Value = 1.029501620854521
$ws.Cells.Item(21, 14).Value = 1.011489064196724

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014822825750009
$ws.Cells.Item(22, 4).Value = 1.022121739783
$ws.Cells.Item(22, 5).Value = 1.0187129889678
$ws.Cells.Item(22, 6).Value = 1.025030145337022
$ws.Cells.Item(22, 9).Value = 1.028737294387917
$ws.Cells.Item(22, 10).Value = 1.021815087359768
$ws.Cells.Item(22, 11).Value = 1.025885304745451
$ws.Cells.Item(22, 12).Value = 1.02249017938402
$ws.Cells.Item(22, 13).Value = 1.0287821645095
$ws.Cells.Item(22, 14).Value = 1.011344544695329

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.015187866023146
$ws.Cells.Item(23, 4).Value = 1.022397456673319
$ws.Cells.Item(23, 5).Value = 1.019056672323724
$ws.Cells.Item(23, 6).Value = 1.025482484223533
$ws.Cells.Item(23, 9).Value = 1.028821734529535
$ws.Cells.Item(23, 10).Value = 1.022044770593631
$ws.Cells.Item(23, 11).Value = 1.026090587436513
$ws.Cells.Item(23, 12).Value = 1.022762885186878
$ws.Cells.Item(23, 13).Value = 1.02916361845198
$ws.Cells.Item(23, 14).Value = 1.011421178653348

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016625134935997
$ws.Cells.Item(24, 4).Value = 1.023482156709369
$ws.Cells.Item(24, 5).Value = 1.020410713673031
$ws.Cells.Item(24, 6).Value = 1.027263282480294
$ws.Cells.Item(24, 9).Value = 1.029150039111204
$ws.Cells.Item(24, 10).Value = 1.022947894769455
$ws.Cells.Item(24, 11).Value = 1.026896490474776
$ws.Cells.Item(24, 12).Value = 1.023836093485474
$ws.Cells.Item(24, 13).Value = 1.030664117976559
$ws.Cells.Item(24, 14).Value = 1.011722403241872

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.018293292980563
$ws.Cells.Item(25, 4).Value = 1.024739263739635
$ws.Cells.Item(25, 5).Value = 1.021984038987519
$ws.Cells.Item(25, 6).Value = 1.029329798974343
$ws.Cells.Item(25, 9).Value = 1.029522455592932
$ws.Cells.Item(25, 10).Value = 1.023993594696244
$ws.Cells.Item(25, 11).Value = 1.028565274147346
$ws.Cells.Item(25, 12).Value = 1.025080630478191
$ws.Cells.Item(25, 13).Value = 1.032402802066844
$ws.Cells.Item(25, 14).Value = 1.012070965829686
